# PARGT_Manual_Mac.docx edit:
#  1. Remove the stray "_GoBack" bookmark that currently sits right after
#     "...save it in a directory" (before the trailing period).
#  2. Change "blastpgp" to "blastpgp_universal-macosx" in the
#     'sudo chmod 755 "blastpgp"' sentence, and leave a "_GoBack" bookmark
#     immediately after the new text (before the closing curly quote).

$d = $word.ActiveDocument

# Unicode curly quotes used throughout the document.
$lq = [char]0x201C   # “
$rq = [char]0x201D   # ”

# --- Step 1: drop the old "_GoBack" bookmark -----------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Step 2: rename blastpgp -> blastpgp_universal-macosx ----------------
$rng = $d.Content
$found = $rng.Find.Execute("blastpgp", $true, $false, $false, $false, $false, `
                            $true, 1, $false, $null, 0)
if ($found) {
    $rng.InsertAfter("_universal-macosx")
}

# --- Step 3: re-add "_GoBack" right after the new text --------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("blastpgp_universal-macosx", $true, $false, $false, `
                              $false, $false, $true, 1, $false, $null, 0)
if ($found2) {
    $rng2.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $rng2)
}
